# forecast_series_OPA_yoy_AR2_50_9_full.xlsx
# Bugfixed evaluation and simulated rt_data for components:
# a new leading observation (date 39400 / FY2007) is inserted at row 2,
# pushing every later row down by one, and the y_0_forecast / y_1_forecast
# columns are recomputed for the whole series (new row 53 added at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new final row (53) the same cell formatting as the other data rows
# (column A keeps the date-formatted style) before filling in values.
$ws.Range("A3:E3").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# date_of_forecast, y_0, y_0_forecast, y_1, y_1_forecast
$data = @(
  (39400, 2007, 5.896808312953783, 2008, $null),
  (39583, 2008, $null, 2009, $null),
  (39765, 2008, 7.441962824572235, 2009, $null),
  (39948, 2009, $null, 2010, $null),
  (40130, 2009, 6.277541464866987, 2010, $null),
  (40310, 2010, 5.907218141265402, 2011, 5.917486466529609),
  (40494, 2010, 6.535114773304773, 2011, 6.325696408067327),
  (40676, 2011, 6.511263427347003, 2012, 6.422943767670297),
  (40862, 2011, 5.12051970717502, 2012, 4.950888348161886),
  (41044, 2012, 4.325828829470257, 2013, 4.493586323244281),
  (41228, 2012, 3.65682115264816, 2013, 3.982564147794321),
  (41409, 2013, 2.964652118442834, 2014, 3.83627393798931),
  (41592, 2013, 2.943878639034381, 2014, 4.334309403335435),
  (41774, 2014, 3.171852776411788, 2015, 4.030605385534614),
  (41957, 2014, 1.172679597477866, 2015, 2.644356903452572),
  (42137, 2015, 1.773712379859993, 2016, 3.484530515673856),
  (42321, 2015, 2.961845079861303, 2016, 3.383932287548697),
  (42503, 2016, 2.533350906619081, 2017, 3.524103740130435),
  (42689, 2016, 2.508469427909898, 2017, 3.355044026998955),
  (42867, 2017, 2.661040979345697, 2018, 3.567108445582057),
  (43053, 2017, 3.523703831572056, 2018, 3.74984170812418),
  (43145, 2018, 2.788213251109939, 2019, 3.595252063027843),
  (43235, 2018, 3.150198973767537, 2019, 3.699072253610103),
  (43326, 2018, 1.334931695392405, 2019, 2.689909849380556),
  (43418, 2018, 1.178605266817589, 2019, 2.186196327763934),
  (43510, 2019, 1.80511617406458, 2020, 3.037278170871094),
  (43600, 2019, 0.4641929091049102, 2020, 2.550259844884462),
  (43691, 2019, 2.9927258084951, 2020, 2.807906319450781),
  (43783, 2019, 3.047037961814492, 2020, 2.880436144359444),
  (43875, 2020, 3.113514644866355, 2021, 2.798268189979214),
  (43966, 2020, 2.585454129751663, 2021, 2.671828487424377),
  (44068, 2020, -0.4891791466461126, 2021, 1.156476476708135),
  (44159, 2020, -0.2228847697281378, 2021, 1.982741503124119),
  (44251, 2021, 2.240953541724267, 2022, 2.598498189609066),
  (44341, 2021, -0.4532848472497908, 2022, 2.066462658785673),
  (44432, 2021, -0.7941560676977599, 2022, 2.320050994894562),
  (44525, 2021, -1.165854108406617, 2022, 2.782217648649521),
  (44617, 2022, 3.392010093835562, 2023, 2.595837839692172),
  (44706, 2022, 0.9477102747197819, 2023, 1.83067479293082),
  (44798, 2022, 2.34069710769782, 2023, 2.116186509693896),
  (44890, 2022, 2.501311189006916, 2023, 2.985901060752827),
  (44981, 2023, 2.735256324140778, 2024, 2.169811364059249),
  (45071, 2023, 2.222852754198135, 2024, 1.934107558751452),
  (45163, 2023, 0.9259311313598806, 2024, 1.280651803859989),
  (45254, 2023, 0.6753076481029074, 2024, 0.7957830962485257),
  (45345, 2024, 1.324987171138314, 2025, 1.882659757536698),
  (45436, 2024, -0.007094633234694392, 2025, 1.444584248586422),
  (45534, 2024, 1.935025917091848, 2025, 2.004076350201744),
  (45618, 2024, 2.039329803030121, 2025, 2.510359031091491),
  (45713, 2025, 1.751699155751707, 2026, 1.973546706924467),
  (45800, 2025, 3.078872076370009, 2026, 2.279508996785351),
  (45891, 2025, 2.43119486791763, 2026, 2.421949074001883)
)

$row = 2
foreach ($values in $data) {
    $cols = @("A", "B", "C", "D", "E")
    for ($i = 0; $i -lt 5; $i++) {
        $cell = $ws.Range($cols[$i] + $row)
        if ($null -eq $values[$i]) {
            $cell.ClearContents()
        } else {
            $cell.Value = $values[$i]
        }
    }
    $row++
}